# Edit "LIST" sheet per commit: numeric stepID version.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

# A2 now holds "AD.SEC.001.FON.01" (previously "AD.SEC.001.FON.02")
$ws.Range("A2").Value = "AD.SEC.001.FON.01"

# A3 now holds "AD.SEC.014.FON.01" (previously "AD.SEC.001.FON.01"),
# and loses the style that made it bold/shaded (now plain/no style).
$ws.Range("A3").Value = "AD.SEC.014.FON.01"
$ws.Range("A3").Style = "Normal"

# A4, A5, A6 are cleared entirely (their cells removed from the sheet).
$ws.Range("A4").Clear()
$ws.Range("A5").Clear()
$ws.Range("A6").Clear()

# New cell H10 gets a new shared string value "MP.CPT.001.CRE"
# (text-format style, matching the other H column entries).
$ws.Range("H10").Value = "MP.CPT.001.CRE"
$ws.Range("H10").NumberFormat = "@"

# New empty styled cell B13 (text-format style, matching B3).
$ws.Range("B13").NumberFormat = "@"

# Update the active selection to C7.
$ws.Range("C7").Select()
